# Generate Report for Handback
# Adds a newly-handed-back file (d19dcf06-dbd6-4f9c-bb07-c47f72314a55) as a
# fresh row on every sheet, and updates the in-flight file's identity from
# 0a24f56a-195f-4297-a987-da646670bcb0 to 12ca941b-7f81-4441-8d13-5e66f3f70501
# now that it has progressed (new xliff hash + timestamps).

$wb = $excel.ActiveWorkbook

$oldFile = "0a24f56a-195f-4297-a987-da646670bcb0"
$renamedFile = "12ca941b-7f81-4441-8d13-5e66f3f70501"
$newFile = "d19dcf06-dbd6-4f9c-bb07-c47f72314a55"

$dateTimeFormat = "yyyy-mm-dd HH:mm:ss"

# =========================================================================
# Sheet "Overview"
# =========================================================================
$wsO = $wb.Worksheets.Item("Overview")
$loO = $wsO.ListObjects.Item(1)
$loO.ListRows.Add()

# Row 2 (renamed file) text updates
$wsO.Range("A2").Value = "$renamedFile.md"
$wsO.Range("C2").Value = ".md"
$wsO.Range("E2").Value = "Handed back: in sync with en-US"
$wsO.Range("F2").Value = "Handed back: in sync with en-US"
$wsO.Range("G2").Value = "2016-09-01 19:11:56"
$wsO.Range("G2").NumberFormat = $dateTimeFormat

# Row 3 (brand-new file)
$wsO.Range("A3").Value = "$newFile.md"
$wsO.Range("C3").Value = ".md"
$wsO.Range("E3").Value = "Handed back: in sync with en-US"
$wsO.Range("F3").Value = "Handed back: in sync with en-US"
$wsO.Range("G3").Value = "2016-09-01 19:11:56"
$wsO.Range("G3").NumberFormat = $dateTimeFormat

# Hyperlinks (column B): refresh the existing one, add the new one
$wsO.Range("B2").Hyperlinks.Delete()
$wsO.Hyperlinks.Add($wsO.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1fecb61ca8a0badef573aba574e902cacc70e0e4/e2e/$renamedFile.md", "", "", "e2e\$renamedFile.md")
$wsO.Hyperlinks.Add($wsO.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1fecb61ca8a0badef573aba574e902cacc70e0e4/e2e/$newFile.md", "", "", "e2e\$newFile.md")

# =========================================================================
# Sheet "zh-cn"
# =========================================================================
$wsZ = $wb.Worksheets.Item("zh-cn")
$loZ = $wsZ.ListObjects.Item(1)
$loZ.ListRows.Add()

# Row 2 (renamed file, handback xliff updated)
$wsZ.Range("A2").Value = "$renamedFile.md"
$wsZ.Range("B2").Value = ".md"
$wsZ.Range("C2").Value = "Handed back: in sync with en-US"
$wsZ.Range("D2").Value = "e2e"
$wsZ.Range("E2").Value = "ht"
$wsZ.Range("F2").Value = "False"
$wsZ.Range("G2").Value = "$renamedFile.91c43bcfdb4ab9074f949a6e65273ec550d206fa.zh-cn.xlf"
$wsZ.Range("H2").Value = "2016-09-01 19:11:50"
$wsZ.Range("H2").NumberFormat = $dateTimeFormat
$wsZ.Range("I2").Value = "$renamedFile.md"
$wsZ.Range("J2").Value = "$renamedFile.91c43bcfdb4ab9074f949a6e65273ec550d206fa.zh-cn.xlf"
$wsZ.Range("K2").Value = "2016-09-01 19:12:22"
$wsZ.Range("K2").NumberFormat = $dateTimeFormat
$wsZ.Range("L2").Value = ""
$wsZ.Range("M2").Value = "True"
$wsZ.Range("N2").Value = ""
$wsZ.Range("O2").Value = "False"
$wsZ.Range("P2").Value = ""

# Row 3 (brand-new file)
$wsZ.Range("A3").Value = "$newFile.md"
$wsZ.Range("B3").Value = ".md"
$wsZ.Range("C3").Value = "Handed back: in sync with en-US"
$wsZ.Range("D3").Value = "e2e"
$wsZ.Range("E3").Value = "ht"
$wsZ.Range("F3").Value = "True"
$wsZ.Range("G3").Value = "$newFile.c2d2aba0b8dc6ae2fdf0f25376e671d7a73970e5.zh-cn.xlf"
$wsZ.Range("H3").Value = "2016-09-01 19:11:50"
$wsZ.Range("H3").NumberFormat = $dateTimeFormat
$wsZ.Range("I3").Value = "$newFile.md"
$wsZ.Range("J3").Value = "$newFile.c2d2aba0b8dc6ae2fdf0f25376e671d7a73970e5.zh-cn.xlf"
$wsZ.Range("K3").Value = "2016-09-01 19:12:22"
$wsZ.Range("K3").NumberFormat = $dateTimeFormat
$wsZ.Range("L3").Value = ""
$wsZ.Range("M3").Value = "True"
$wsZ.Range("N3").Value = ""
$wsZ.Range("O3").Value = "False"
$wsZ.Range("P3").Value = ""

# Hyperlinks (columns A & I): refresh existing, add new
$wsZ.Range("A2").Hyperlinks.Delete()
$wsZ.Range("I2").Hyperlinks.Delete()
$wsZ.Hyperlinks.Add($wsZ.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1fecb61ca8a0badef573aba574e902cacc70e0e4/e2e/$renamedFile.md", "", "", "$renamedFile.md")
$wsZ.Hyperlinks.Add($wsZ.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/2fed95f712d213262fa1ed4764bc89b7f67d7336/e2e/$renamedFile.md", "", "", "$renamedFile.md")
$wsZ.Hyperlinks.Add($wsZ.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1fecb61ca8a0badef573aba574e902cacc70e0e4/e2e/$newFile.md", "", "", "$newFile.md")
$wsZ.Hyperlinks.Add($wsZ.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/2fed95f712d213262fa1ed4764bc89b7f67d7336/e2e/$newFile.md", "", "", "$newFile.md")

# =========================================================================
# Sheet "de-de"
# =========================================================================
$wsD = $wb.Worksheets.Item("de-de")
$loD = $wsD.ListObjects.Item(1)
$loD.ListRows.Add()

# Row 2 (renamed file, handback xliff updated)
$wsD.Range("A2").Value = "$renamedFile.md"
$wsD.Range("B2").Value = ".md"
$wsD.Range("C2").Value = "Handed back: in sync with en-US"
$wsD.Range("D2").Value = "e2e"
$wsD.Range("E2").Value = "ht"
$wsD.Range("F2").Value = "False"
$wsD.Range("G2").Value = "$renamedFile.91c43bcfdb4ab9074f949a6e65273ec550d206fa.de-de.xlf"
$wsD.Range("H2").Value = "2016-09-01 19:11:56"
$wsD.Range("H2").NumberFormat = $dateTimeFormat
$wsD.Range("I2").Value = "$renamedFile.md"
$wsD.Range("J2").Value = "$renamedFile.91c43bcfdb4ab9074f949a6e65273ec550d206fa.de-de.xlf"
$wsD.Range("K2").Value = "2016-09-01 19:12:30"
$wsD.Range("K2").NumberFormat = $dateTimeFormat
$wsD.Range("L2").Value = ""
$wsD.Range("M2").Value = "True"
$wsD.Range("N2").Value = ""
$wsD.Range("O2").Value = "False"
$wsD.Range("P2").Value = ""

# Row 3 (brand-new file)
$wsD.Range("A3").Value = "$newFile.md"
$wsD.Range("B3").Value = ".md"
$wsD.Range("C3").Value = "Handed back: in sync with en-US"
$wsD.Range("D3").Value = "e2e"
$wsD.Range("E3").Value = "ht"
$wsD.Range("F3").Value = "True"
$wsD.Range("G3").Value = "$newFile.c2d2aba0b8dc6ae2fdf0f25376e671d7a73970e5.de-de.xlf"
$wsD.Range("H3").Value = "2016-09-01 19:11:56"
$wsD.Range("H3").NumberFormat = $dateTimeFormat
$wsD.Range("I3").Value = "$newFile.md"
$wsD.Range("J3").Value = "$newFile.c2d2aba0b8dc6ae2fdf0f25376e671d7a73970e5.de-de.xlf"
$wsD.Range("K3").Value = "2016-09-01 19:12:30"
$wsD.Range("K3").NumberFormat = $dateTimeFormat
$wsD.Range("L3").Value = ""
$wsD.Range("M3").Value = "True"
$wsD.Range("N3").Value = ""
$wsD.Range("O3").Value = "False"
$wsD.Range("P3").Value = ""

# Hyperlinks (columns A & I): refresh existing, add new
$wsD.Range("A2").Hyperlinks.Delete()
$wsD.Range("I2").Hyperlinks.Delete()
$wsD.Hyperlinks.Add($wsD.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1fecb61ca8a0badef573aba574e902cacc70e0e4/e2e/$renamedFile.md", "", "", "$renamedFile.md")
$wsD.Hyperlinks.Add($wsD.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/145ae36dad905ad30435028fe029e4ecc95d66eb/e2e/$renamedFile.md", "", "", "$renamedFile.md")
$wsD.Hyperlinks.Add($wsD.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1fecb61ca8a0badef573aba574e902cacc70e0e4/e2e/$newFile.md", "", "", "$newFile.md")
$wsD.Hyperlinks.Add($wsD.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/145ae36dad905ad30435028fe029e4ecc95d66eb/e2e/$newFile.md", "", "", "$newFile.md")

Write-Host "Handback status report updated."
